# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the aggregated "全部类型" sheet, matching the regenerated site data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 959
$ws.Range("F7").Value = 3037
$ws.Range("F8").Value = 1923
$ws.Range("F13").Value = 281
$ws.Range("F18").Value = 9663
$ws.Range("F20").Value = 9
$ws.Range("F22").Value = 7635
$ws.Range("F23").Value = 12175
$ws.Range("F26").Value = 250
$ws.Range("F29").Value = 2741
$ws.Range("F30").Value = 245
$ws.Range("F32").Value = 2743
$ws.Range("F33").Value = 1202
$ws.Range("F35").Value = 62
$ws.Range("F37").Value = 4567
$ws.Range("F38").Value = 1175
$ws.Range("F42").Value = 587

# --- Sheet "全部类型" (all types, aggregated) ---
$ws2 = $wb.Worksheets.Item("全部类型")
$ws2.Range("F7").Value = 959
$ws2.Range("F10").Value = 3037
$ws2.Range("F11").Value = 1923
$ws2.Range("F17").Value = 281
$ws2.Range("F22").Value = 9663
$ws2.Range("F25").Value = 7635
$ws2.Range("F26").Value = 12175
$ws2.Range("F28").Value = 250
$ws2.Range("F32").Value = 2741
$ws2.Range("F35").Value = 245
$ws2.Range("F36").Value = 62
$ws2.Range("F38").Value = 4567
$ws2.Range("F46").Value = 587
